# Metodos create y delete
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "Reunión" default type with "Por defecto" for the event rows
$ws.Range("D2").Value = "Por defecto"
$ws.Range("D3").Value = "Por defecto"
$ws.Range("D4").Value = "Por defecto"

# Update the active selection to D2 as in the final file
$ws.Range("D2").Select()
